$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.950.49"
$ws.Range("E2").Value = "  -1.13%  "

$ws.Range("D3").Value = "3.395.18"
$ws.Range("E3").Value = "  -1.62%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'572.86"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").Value = "'142.41"
$ws.Range("E6").Value = "  -2.31%  "

$ws.Range("D7").Value = "3.395.18"
$ws.Range("E7").Value = "  -1.66%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").Value = "'7.58"
$ws.Range("E10").Value = "  -1.65%  "

$ws.Range("E11").Value = "  -2.43%  "

$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").Value = "3.973.39"
$ws.Range("E13").Value = "  -1.62%  "

$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").Value = "'28.06"
$ws.Range("E15").Value = "  -1.51%  "

$ws.Range("E16").Value = "  -1.37%  "

$ws.Range("D17").Value = "3.398.95"
$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("D18").Value = "61.000.67"
$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("E19").Value = "  -3.90%  "

$ws.Range("D20").Value = "'13.85"
$ws.Range("E20").Value = "  -3.43%  "

$ws.Range("D21").Value = "'8.97"
$ws.Range("E21").Value = "  -4.87%  "

$ws.Range("D22").Value = "'382.71"
$ws.Range("E22").Value = "  -5.08%  "

$ws.Range("E23").Value = "  -1.72%  "

$ws.Range("D24").Value = "'74.31"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("E26").Value = "  -5.40%  "

$ws.Range("D27").Value = "3.530.61"
$ws.Range("E27").Value = "  -1.61%  "

$ws.Range("D28").Value = "'0.180"
$ws.Range("E28").Value = "  -1.86%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").Value = "'7.37"
$ws.Range("E30").Value = "  -3.33%  "

$ws.Range("D31").Value = "'8.03"
$ws.Range("E31").Value = "  -2.78%  "

$ws.Range("E32").Value = "  -1.38%  "

$ws.Range("D33").Value = "'1.42"
$ws.Range("E33").Value = "  -2.56%  "

$ws.Range("D35").Value = "'23.46"
$ws.Range("E35").Value = "  -2.08%  "

$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("D37").Value = "'167.74"
$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").Value = "3.424.08"
$ws.Range("E38").Value = "  -1.51%  "

$ws.Range("E39").Value = "  -2.51%  "

$ws.Range("D40").Value = "'1.49"
$ws.Range("E40").Value = "  -5.21%  "

$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").Value = "'27.21"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  -2.53%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("E45").Value = "  -2.12%  "

$ws.Range("E46").Value = "  -4.45%  "

$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("D48").Value = "2.489.33"
$ws.Range("E48").Value = "  -4.59%  "

$ws.Range("D49").Value = "'6.82"
$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("D50").Value = "'23.05"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("E51").Value = "  +0.79%  "
